$wb = $excel.ActiveWorkbook

# This script applies a data refresh to the cached market-price / profit
# columns (H..N) produced by the scheduled Universalis price-fetch runner.
# Values are plain numeric literals (no formulas in the source data).

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 6063109.5
$ws.Range("I106").Value = 6063109.5
$ws.Range("K106").Value = 6063109.5
$ws.Range("M106").Value = -6062478.5
$ws.Range("H112").Value = 1518746.4
$ws.Range("J112").Value = 1758174.8
$ws.Range("L112").Value = 5274524.4
$ws.Range("N112").Value = -5276740.4
$ws.Range("H125").Value = 3696.2856
$ws.Range("I125").Value = 687.6667
$ws.Range("J125").Value = 5952.75
$ws.Range("K125").Value = 6189.0003
$ws.Range("L125").Value = 53574.75
$ws.Range("M125").Value = -3729.0003
$ws.Range("N125").Value = -58494.75
$ws.Range("H127").Value = 829.2
$ws.Range("I127").Value = 829.2
$ws.Range("K127").Value = 2487.6
$ws.Range("M127").Value = 2472.4
$ws.Range("H129").Value = 1340.1666
$ws.Range("I129").Value = 812.3
$ws.Range("K129").Value = 2436.9
$ws.Range("M129").Value = 2563.1
$ws.Range("H137").Value = 1738156.9
$ws.Range("I137").Value = 1837.1111
$ws.Range("K137").Value = 5511.3333
$ws.Range("M137").Value = -2961.3333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26368272
$ws.Range("I32").Value = 29612014
$ws.Range("J32").Value = 8933160
$ws.Range("K32").Value = 29612014
$ws.Range("L32").Value = 8933160
$ws.Range("M32").Value = -29611727
$ws.Range("N32").Value = -8933734
$ws.Range("H61").Value = 3294
$ws.Range("I61").Value = 2896
$ws.Range("J61").Value = 3661.3845
$ws.Range("K61").Value = 2896
$ws.Range("L61").Value = 3661.3845
$ws.Range("M61").Value = -2684
$ws.Range("N61").Value = -4085.3845
$ws.Range("H63").Value = 4829.231
$ws.Range("J63").Value = 4829.231
$ws.Range("L63").Value = 4829.231
$ws.Range("N63").Value = -6201.231
$ws.Range("H66").Value = 4829.231
$ws.Range("J66").Value = 4829.231
$ws.Range("L66").Value = 24146.155
$ws.Range("N66").Value = -31010.155
$ws.Range("H74").Value = 2883.25
$ws.Range("I74").Value = 2343.1853
$ws.Range("J74").Value = 5799.6
$ws.Range("K74").Value = 2343.1853
$ws.Range("L74").Value = 5799.6
$ws.Range("M74").Value = -1469.1853
$ws.Range("N74").Value = -7547.6
$ws.Range("H77").Value = 2883.25
$ws.Range("I77").Value = 2343.1853
$ws.Range("J77").Value = 5799.6
$ws.Range("K77").Value = 11715.9265
$ws.Range("L77").Value = 28998
$ws.Range("M77").Value = -7347.926500000001
$ws.Range("N77").Value = -37734
$ws.Range("H88").Value = 1741.3334
$ws.Range("I88").Value = 2299
$ws.Range("J88").Value = 1462.5
$ws.Range("K88").Value = 2299
$ws.Range("L88").Value = 1462.5
$ws.Range("M88").Value = -1893
$ws.Range("N88").Value = -2274.5
$ws.Range("H91").Value = 1741.3334
$ws.Range("I91").Value = 2299
$ws.Range("J91").Value = 1462.5
$ws.Range("K91").Value = 2299
$ws.Range("L91").Value = 1462.5
$ws.Range("M91").Value = -895
$ws.Range("N91").Value = -4270.5
$ws.Range("H102").Value = 2388
$ws.Range("I102").Value = 1984.0834
$ws.Range("K102").Value = 1984.0834
$ws.Range("M102").Value = -362.0834
$ws.Range("H132").Value = 4022.1538
$ws.Range("I132").Value = 3120.0356
$ws.Range("K132").Value = 9360.106800000001
$ws.Range("M132").Value = -6830.106800000001
$ws.Range("H136").Value = 3294
$ws.Range("I136").Value = 2896
$ws.Range("J136").Value = 3661.3845
$ws.Range("K136").Value = 8688
$ws.Range("L136").Value = 10984.1535
$ws.Range("M136").Value = -6138
$ws.Range("N136").Value = -16084.1535

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 170666.67
$ws.Range("I20").Value = 500000
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 500000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -499753
$ws.Range("N20").Value = -6494
$ws.Range("H86").Value = 2999
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2999
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2999
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5245
$ws.Range("H89").Value = 2999
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2999
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 14995
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -26227
$ws.Range("H105").Value = 2670.8462
$ws.Range("I105").Value = 2311.7058
$ws.Range("K105").Value = 2311.7058
$ws.Range("M105").Value = -564.7058000000002
$ws.Range("H110").Value = 55957
$ws.Range("J110").Value = 55957
$ws.Range("L110").Value = 55957
$ws.Range("N110").Value = -64137
$ws.Range("H112").Value = 152000
$ws.Range("J112").Value = 152000
$ws.Range("L112").Value = 152000
$ws.Range("N112").Value = -154954

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 123000
$ws.Range("J20").Value = 123000
$ws.Range("L20").Value = 123000
$ws.Range("N20").Value = -123472
$ws.Range("H30").Value = 123000
$ws.Range("J30").Value = 123000
$ws.Range("L30").Value = 123000
$ws.Range("N30").Value = -123182
$ws.Range("H31").Value = 5393.3447
$ws.Range("J31").Value = 9916.666999999999
$ws.Range("L31").Value = 9916.666999999999
$ws.Range("N31").Value = -10506.667
$ws.Range("H34").Value = 5393.3447
$ws.Range("J34").Value = 9916.666999999999
$ws.Range("L34").Value = 9916.666999999999
$ws.Range("N34").Value = -10320.667
$ws.Range("H58").Value = 2701.6736
$ws.Range("I58").Value = 2487.5854
$ws.Range("J58").Value = 3798.875
$ws.Range("K58").Value = 2487.5854
$ws.Range("L58").Value = 3798.875
$ws.Range("M58").Value = -2284.5854
$ws.Range("N58").Value = -4204.875
$ws.Range("H112").Value = 99900
$ws.Range("J112").Value = 99900
$ws.Range("L112").Value = 99900
$ws.Range("N112").Value = -102854
$ws.Range("H128").Value = 123000
$ws.Range("J128").Value = 123000
$ws.Range("L128").Value = 123000
$ws.Range("N128").Value = -132960
$ws.Range("H136").Value = 2701.6736
$ws.Range("I136").Value = 2487.5854
$ws.Range("J136").Value = 3798.875
$ws.Range("K136").Value = 7462.7562
$ws.Range("L136").Value = 11396.625
$ws.Range("M136").Value = -4912.7562
$ws.Range("N136").Value = -16496.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3282.7144
$ws.Range("J75").Value = 3163.1667
$ws.Range("L75").Value = 9489.500100000001
$ws.Range("N75").Value = -11485.5001
$ws.Range("H78").Value = 3282.7144
$ws.Range("J78").Value = 3163.1667
$ws.Range("L78").Value = 28468.5003
$ws.Range("N78").Value = -38452.5003
$ws.Range("H113").Value = 1270.25
$ws.Range("I113").Value = 567
$ws.Range("J113").Value = 1504.6666
$ws.Range("K113").Value = 1701
$ws.Range("L113").Value = 4513.9998
$ws.Range("M113").Value = 469
$ws.Range("N113").Value = -8853.9998
$ws.Range("H117").Value = 3537
$ws.Range("I117").Value = 2163.75
$ws.Range("J117").Value = 4147.3335
$ws.Range("K117").Value = 6491.25
$ws.Range("L117").Value = 12442.0005
$ws.Range("M117").Value = -3049.25
$ws.Range("N117").Value = -19326.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 69947.5
$ws.Range("J64").Value = 69947.5
$ws.Range("L64").Value = 69947.5
$ws.Range("N64").Value = -70443.5
$ws.Range("H67").Value = 69947.5
$ws.Range("J67").Value = 69947.5
$ws.Range("L67").Value = 69947.5
$ws.Range("N67").Value = -71663.5
$ws.Range("H70").Value = 83434.375
$ws.Range("I70").Value = 107746
$ws.Range("J70").Value = 10499.5
$ws.Range("K70").Value = 107746
$ws.Range("L70").Value = 10499.5
$ws.Range("M70").Value = -107476
$ws.Range("N70").Value = -11039.5
$ws.Range("H73").Value = 83434.375
$ws.Range("I73").Value = 107746
$ws.Range("J73").Value = 10499.5
$ws.Range("K73").Value = 107746
$ws.Range("L73").Value = 10499.5
$ws.Range("M73").Value = -106810
$ws.Range("N73").Value = -12371.5
$ws.Range("H126").Value = 2619.625
$ws.Range("I126").Value = 2493.8572
$ws.Range("K126").Value = 7481.571599999999
$ws.Range("M126").Value = -5011.571599999999
$ws.Range("H128").Value = 112330
$ws.Range("J128").Value = 112330
$ws.Range("L128").Value = 112330
$ws.Range("N128").Value = -122290
$ws.Range("H132").Value = 4272.1816
$ws.Range("I132").Value = 3374.375
$ws.Range("K132").Value = 10123.125
$ws.Range("M132").Value = -7593.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 32336.637
$ws.Range("I122").Value = 33967.11
$ws.Range("J122").Value = 24999.5
$ws.Range("K122").Value = 101901.33
$ws.Range("L122").Value = 74998.5
$ws.Range("M122").Value = -99451.33
$ws.Range("N122").Value = -79898.5
$ws.Range("H131").Value = 107496.336
$ws.Range("J131").Value = 107496.336
$ws.Range("L131").Value = 107496.336
$ws.Range("N131").Value = -117576.336
$ws.Range("H132").Value = 3715.4167
$ws.Range("I132").Value = 3658.6
$ws.Range("K132").Value = 10975.8
$ws.Range("M132").Value = -8445.799999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 6009
$ws.Range("I30").Value = 6009
$ws.Range("K30").Value = 6009
$ws.Range("M30").Value = -5902
$ws.Range("H37").Value = 27507
$ws.Range("I37").Value = 20763
$ws.Range("J37").Value = 40995
$ws.Range("K37").Value = 20763
$ws.Range("L37").Value = 40995
$ws.Range("M37").Value = -20560
$ws.Range("N37").Value = -41401
$ws.Range("H122").Value = 166672910
$ws.Range("I122").Value = 166672910
$ws.Range("K122").Value = 500018730
$ws.Range("M122").Value = -500016280
$ws.Range("H132").Value = 3505.0625
$ws.Range("I132").Value = 3506.3076
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 10518.9228
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -7988.9228
$ws.Range("N132").Value = -15559.0001
